$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.595.99"
$ws.Range("E2").Value = "  +5.85%  "
$ws.Range("D3").Value = "2.741.60"
$ws.Range("E3").Value = "  +4.82%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'116.12"
$ws.Range("E5").Value = "  +5.69%  "
$ws.Range("D6").Value = "'333.82"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +5.20%  "
$ws.Range("D10").Value = "'41.36"
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("E11").Value = "  +5.45%  "
$ws.Range("D12").Value = "'20.11"
$ws.Range("E12").Value = "  +2.37%  "
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("E14").Value = "  +4.47%  "
$ws.Range("D15").Value = "3.171.56"
$ws.Range("E15").Value = "  +4.90%  "
$ws.Range("D16").Value = "2.772.40"
$ws.Range("E16").Value = "  +5.65%  "
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "51.503.03"
$ws.Range("E18").Value = "  +5.78%  "
$ws.Range("E19").Value = "  +5.44%  "
$ws.Range("D20").Value = "'13.43"
$ws.Range("E20").Value = "  +4.86%  "
$ws.Range("D21").Value = "'6.82"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "0.0₃0970"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").Value = "'279.04"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").Value = "'69.15"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "'2.65"
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("D26").Value = "'26.67"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'0.140"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "'34.79"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'49.92"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "'0.0817"
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "'18.91"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "'2.07"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").Value = "'3.17"
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("D40").Value = "'127.52"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("E41").Value = "  +8.71%  "
$ws.Range("D42").Value = "'22.93"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("D44").Value = "'2.28"
$ws.Range("E44").Value = "  +7.19%  "
$ws.Range("E45").Value = "  +12.29%  "
$ws.Range("D46").Value = "2.085.76"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("E49").Value = "  +7.01%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'59.59"
$ws.Range("E51").Value = "  +2.07%  "
